# "0.1 version - It is working at a very basic level"
#
#  - Rename the "Reference" sheet to "vMixConfig". Excel automatically
#    rewrites every formula/defined name that pointed at "Reference"
#    (INDEX/MATCH lookups on the Plan sheet, the VMix/mediaNumber/mediaType/
#    shortTitle/etc. defined names, ...) to point at "vMixConfig" instead.
#  - Make "vMixConfig" the active/selected sheet (previously "Plan" was the
#    selected tab).

$wb = $excel.ActiveWorkbook

$refSheet = $wb.Worksheets.Item("Reference")

# Renaming cascades through all dependent formulas and defined names.
$refSheet.Name = "vMixConfig"

# Switch the active tab from "Plan" to the newly renamed sheet.
$refSheet.Activate()
$refSheet.Select()

# Best-effort: keep the list-validation on Plan!F2:F47 (sourced from the
# renamed sheet's column T) pointing at the new sheet name too.
try {
    $planSheet = $wb.Worksheets.Item("Plan")
    $planSheet.Range("F2:F47").Validation.Formula1 = "=vMixConfig!`$T:`$T"
} catch {
    # Older/limited validation object models may not expose this; the sheet
    # rename above is the substantive part of the change.
}
